# moneybin: calculate the smoothing factor (%D)
# %D is a simple moving average of %K over a defined smoothing period.

$wb = $excel.ActiveWorkbook
$wsParams = $wb.Worksheets.Item("Parameters")
$wsData = $wb.Worksheets.Item("Data")

# --- Data sheet: add the %D column header ---
$wsData.Range("K1").Value = "%D"

# --- Data sheet: add a new row of stock data (row 17) ---
$wsData.Range("A17").Value = 41407
$wsData.Range("A17").NumberFormat = "M/D/YYYY"
$wsData.Range("B17").Value = 1632.1
$wsData.Range("C17").Value = 1636
$wsData.Range("D17").Value = 1626.74
$wsData.Range("E17").Value = 1633.77
$wsData.Range("F17").Value = 2910600000
$wsData.Range("G17").Value = 1633.77
$wsData.Range("H17").Formula = "=MAX(C4:C17)"
$wsData.Range("I17").Formula = "=MIN(D4:D17)"
$wsData.Range("J17").Formula = "=(E17 - I17) / (H17 - I17) * 100"

# %D: simple moving average of %K (J column) over the smoothing period
$wsData.Range("K17").Formula = "=AVERAGE(J15:J17)"

# Leave the Data sheet's cursor on the new %D column before switching away
$wsData.Activate() | Out-Null
$wsData.Range("K2").Select() | Out-Null

# --- Restore the active sheet/tab to Parameters (with its prior selection) ---
$wsParams.Activate() | Out-Null
$wsParams.Range("C6").Select() | Out-Null
